$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row, derived from the target diff (rows re-sorted chronologically
# by Fecha, with Volumen / Precio columns following their row).
$data = @{
    2 = @{ D = 44508; J = 40; K = 10000; L = 10000; M = 10000; P = 667 }
    3 = @{ D = 44756; J = 80; K = 20000; L = 20000; M = 20000; P = 1333 }
    4 = @{ D = 44525; J = 40; K = 8000;  L = 8000;  M = 8000;  P = 533 }
    5 = @{ D = 44755; J = 50; K = 20000; L = 20000; M = 20000; P = 1333 }
    6 = @{ D = 44749; J = 50; K = 20000; L = 20000; M = 20000; P = 1333 }
    7 = @{ D = 44518; J = 50; K = 10000; L = 10000; M = 10000; P = 667 }
    8 = @{ D = 44757; J = 30; K = 20000; L = 20000; M = 20000; P = 1333 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("P$row").Value = $vals.P
}
